$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: change B43 from text "2" to numeric 2
$ws.Range("B43").Value = 2

# Add new row 44 with annotation data
$ws.Range("A44").Value = "Sunsi Wu"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "2"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "why"
$ws.Range("D44").Value = "FBK"
$ws.Range("E44").Value = "EXP"
$ws.Range("F44").Value = "77ff87fb-cfc5-44ac-a4b7-cb33b05fed6f"
$ws.Range("G44").Value = "ByQpn1ZA-_annotated.xlsx"
$ws.Range("H44").Value = "If we know the regularization is fundamentally and mathematically wrong, why do we investigate its performance?"
